# Registree stats backup on Tue 20 Apr 2021 08:45:05 SAST
#
# - Refresh the "as of" timestamp on every sheet title (19/04/2021 10:15 -> 20/04/2021 08:45)
# - Insert a newly-registered attendee (Strydom, Melinda Lee / Alberton / 410E) into the
#   MD410 Attendance roll at row 183 (alphabetical order), pushing the rest of the list down
# - Bump the "Number of attendees" footer count from 232 to 233

$wb = $excel.ActiveWorkbook

# --- MD410 Attendance ---------------------------------------------------
$ws = $wb.Worksheets.Item("MD410 Attendance")
$ws.Range("A1").Value = "MD410 Registrees as of 20/04/2021 08:45"

# Insert a new data row at 183, shifting the existing rows (old 183.."t Hart"/totals) down one.
$ws.Rows.Item(183).Insert()
$ws.Rows.Item(183).RowHeight = 25
$ws.Range("A183:F183").Borders.LineStyle = 1

$ws.Range("A183").Value = "Strydom"
$ws.Range("B183").Value = "Melinda Lee"
$ws.Range("C183").Value = "Alberton"
$ws.Range("D183").Value = "No"
$ws.Range("E183").Value = "No"
$ws.Range("F183").Value = "410E"

# Footer "Number of attendees" row moved from 235 -> 236 because of the insert above.
$ws.Range("A236").Value = "Number of attendees: 233"

# --- 410E Attendance ------------------------------------------------------
$ws2 = $wb.Worksheets.Item("410E Attendance")
$ws2.Range("A1").Value = "410E Registrees as of 20/04/2021 08:45"

# --- 410W Attendance ------------------------------------------------------
$ws3 = $wb.Worksheets.Item("410W Attendance")
$ws3.Range("A1").Value = "410W Registrees as of 20/04/2021 08:45"

# --- 410E Voting ------------------------------------------------------
$ws4 = $wb.Worksheets.Item("410E Voting")
$ws4.Range("A1").Value = "410E Voting details as of 20/04/2021 08:45"

# --- 410W Voting ------------------------------------------------------
$ws5 = $wb.Worksheets.Item("410W Voting")
$ws5.Range("A1").Value = "410W Voting details as of 20/04/2021 08:45"
